$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "Gemma-7B-Instruct"
$ws.Range("B8").Value = "0.79 ± 0.41"
$ws.Range("C8").Value = "0.19 ± 0.76"
$ws.Range("D8").Value = "0.25 ± 0.56"
$ws.Range("E8").Value = "0.01 ± 0.01"
$ws.Range("F8").Value = "0.09 ± 0.08"
$ws.Range("G8").Value = "0.01 ± 0.04"
$ws.Range("H8").Value = "0.09 ± 0.07"
$ws.Range("I8").Value = "0.12 ± 0.1"
$ws.Range("J8").Value = "0.75 ± 0.26"
$ws.Range("K8").Value = "0.76 ± 0.26"
$ws.Range("L8").Value = "0.75 ± 0.26"
$ws.Range("M8").Value = "0.77 ± 0.27"
$ws.Range("N8").Value = "0.88 ± 0.3"
$ws.Range("O8").Value = "0.12 ± 0.13"
$ws.Range("P8").Value = "0.57 ± 0.24"
$ws.Range("Q8").Value = "8.64 ± 1.59"
$ws.Range("R8").Value = "0.105 ± 0.00"
$ws.Range("S8").Value = "0.85 ± 0.29"
$ws.Range("T8").Value = "0.73 ± 0.4"
$ws.Range("U8").Value = "2.42 ± 1.38"
$ws.Range("V8").Value = "0.74 ± 0.41"
$ws.Range("W8").Value = "0.85 ± 0.29"
$ws.Range("X8").Value = "1.12 ± 0.45"
